# Updated cryptos list on Thu Nov 21 08:47:13 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "96.840.14"
$ws.Range("E2").Value = "  +3.83%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.144.81"
$ws.Range("E3").Value = "  +0.57%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.35"
$ws.Range("E5").Value = "  +1.59%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "612.96"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.78%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.386"
$ws.Range("E8").Value = "  -2.26%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.142.56"
$ws.Range("E10").Value = "  +0.58%  "

# Row 11 - Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.782"
$ws.Range("E11").Value = "  -4.07%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.36%  "

# Row 13 - WrappedBTC
$ws.Range("D13").Value = "96.657.25"
$ws.Range("E13").Value = "  +4.12%  "

# Row 14 - ShibaInu->Toncoin
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.58"
$ws.Range("E14").Value = "  +2.31%  "

# Row 15 - Toncoin->ShibaInu
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000241"
$ws.Range("E15").Value = "  -2.04%  "

# Row 16 - Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.26"
$ws.Range("E16").Value = "  -1.56%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.721.15"
$ws.Range("E17").Value = "  +0.32%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.143.06"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19 - SuiNetwork
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.56"
$ws.Range("E19").Value = "  -6.14%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "517.50"
$ws.Range("E20").Value = "  +16.88%  "

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.67"
$ws.Range("E21").Value = "  +0.13%  "

# Row 22 - Polkadot
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.70"
$ws.Range("E22").Value = "  -4.20%  "

# Row 23 - PEPE
$ws.Range("E23").Value = "  -5.36%  "

# Row 24 - Uniswap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.87"
$ws.Range("E24").Value = "  -3.36%  "

# Row 25 - NEARProtocol
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.52"
$ws.Range("E25").Value = "  -2.18%  "

# Row 26 - Litecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.79"
$ws.Range("E26").Value = "  +3.33%  "

# Row 27 - Aptos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.69"
$ws.Range("E27").Value = "  -5.63%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "3.302.10"
$ws.Range("E28").Value = "  +0.45%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.06%  "

# Row 30 - Stellar
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.239"
$ws.Range("E30").Value = "  +2.41%  "

# Row 31 - Cronos
$ws.Range("E31").Value = "  -2.54%  "

# Row 32 - Hedera
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.125"
$ws.Range("E32").Value = "  +1.26%  "

# Row 33 - Binance-PegBSC-USD
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.49%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.06"
$ws.Range("E34").Value = "  -1.54%  "

# Row 35 - EthereumClassic
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.81"
$ws.Range("E35").Value = "  +3.46%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -4.54%  "

# Row 37 - RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.36"
$ws.Range("E37").Value = "  -9.17%  "

# Row 38 - PancakeSwap
$ws.Range("E38").Value = "  -0.84%  "

# Row 39 - WhiteBITCoin
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.21"
$ws.Range("E39").Value = "  +0.92%  "

# Row 40 - Bittensor
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "480.79"
$ws.Range("E40").Value = "  +1.47%  "

# Row 41 - PolygonEcosystemToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.439"
$ws.Range("E41").Value = "  +1.73%  "

# Row 42 - Fetch.AI
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.23"
$ws.Range("E42").Value = "  -5.21%  "

# Row 43 - MantraDAO
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.56"
$ws.Range("E43").Value = "  -10.41%  "

# Row 45 - dogwifhat
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("E45").Value = "  -4.80%  "

# Row 46 - Monero
$ws.Range("E46").Value = "  +0.89%  "

# Row 47 - Stacks->ARBITRUM
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.712"
$ws.Range("E47").Value = "  +3.16%  "

# Row 48 - ARBITRUM->Stacks
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.92"
$ws.Range("E48").Value = "  +4.10%  "

# Row 49 - Filecoin
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.50"
$ws.Range("E49").Value = "  +2.37%  "

# Row 50 - OKB
$ws.Range("E50").Value = "  +0.38%  "
